$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing B2 cell text
$ws.Range("B2").Value = "sampleFile2.zvi_filtered1.png"

# Copy the style of A2 (bold/border header-like style) down to A3:A6
$ws.Range("A2").Copy()
$ws.Range("A3:A6").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Data for new rows 3-6
$rows = @(
    @{ Row = 3; A = 1; B = "sampleFile3.zvi_filtered2.png"; C = 0.064; D = 0.511904761904762; E = 0.511904761904762 },
    @{ Row = 4; A = 2; B = "sampleFile4.zvi_filtered3.png"; C = 0.064; D = 0.511904761904762; E = 0.511904761904762 },
    @{ Row = 5; A = 3; B = "sampleFile1.zvi_filtered4.png"; C = 0.064; D = 0.511904761904762; E = 0.511904761904762 },
    @{ Row = 6; A = 4; B = "sampleFile5.zvi_filtered5.png"; C = 0.064; D = 0.511904761904762; E = 0.511904761904762 }
)

foreach ($row in $rows) {
    $r = $row.Row
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
}
